# Generate Report for Handback
# Replace the "78b55644-0d76-42fb-a98a-4b0706e7c488" handoff file with the
# newly generated "230d5acb-504b-4637-9f8e-a5c285449e76" file, and the
# "7b2c538e-3fdb-45ed-be82-d28ae0801cc8" file with
# "ffffd46d8d01-c2cf-4e44-ab34-61cf8c3e76e4", refreshing the associated
# xliff names and generate/handback timestamps across all three sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Range("A2").Value = '230d5acb-504b-4637-9f8e-a5c285449e76.md'
$ws1.Range("B2").Value = 'e2e\230d5acb-504b-4637-9f8e-a5c285449e76.md'
$ws1.Range("G2").Value = '2016-09-04 07:09:09'

$ws1.Range("A3").Value = 'ffffd46d8d01-c2cf-4e44-ab34-61cf8c3e76e4.md'
$ws1.Range("B3").Value = 'e2e\ffffd46d8d01-c2cf-4e44-ab34-61cf8c3e76e4.md'
$ws1.Range("G3").Value = '2016-09-04 07:09:09'

$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d15cafe1ec021c40fea94da5b6f49b1a36cc7d1c/e2e/78b55644-0d76-42fb-a98a-4b0706e7c488.md", "", "", 'e2e\230d5acb-504b-4637-9f8e-a5c285449e76.md')
$ws1.Hyperlinks.Add($ws1.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d15cafe1ec021c40fea94da5b6f49b1a36cc7d1c/e2e/7b2c538e-3fdb-45ed-be82-d28ae0801cc8.md", "", "", 'e2e\ffffd46d8d01-c2cf-4e44-ab34-61cf8c3e76e4.md')

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Range("A2").Value = '230d5acb-504b-4637-9f8e-a5c285449e76.md'
$ws2.Range("G2").Value = '230d5acb-504b-4637-9f8e-a5c285449e76.f54858a962c9264f22fab170293849181df9a130.zh-cn.xlf'
$ws2.Range("H2").Value = '2016-09-04 07:09:00'
$ws2.Range("I2").Value = '230d5acb-504b-4637-9f8e-a5c285449e76.md'
$ws2.Range("J2").Value = '230d5acb-504b-4637-9f8e-a5c285449e76.f54858a962c9264f22fab170293849181df9a130.zh-cn.xlf'
$ws2.Range("K2").Value = '2016-09-04 07:09:28'

$ws2.Range("A3").Value = 'ffffd46d8d01-c2cf-4e44-ab34-61cf8c3e76e4.md'
$ws2.Range("G3").Value = '230d5acb-504b-4637-9f8e-a5c285449e76.f54858a962c9264f22fab170293849181df9a130.zh-cn.xlf'
$ws2.Range("H3").Value = '2016-09-04 07:09:00'
$ws2.Range("I3").Value = 'ffffd46d8d01-c2cf-4e44-ab34-61cf8c3e76e4.md'
$ws2.Range("J3").Value = '230d5acb-504b-4637-9f8e-a5c285449e76.f54858a962c9264f22fab170293849181df9a130.zh-cn.xlf'
$ws2.Range("K3").Value = '2016-09-04 07:09:28'

$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d15cafe1ec021c40fea94da5b6f49b1a36cc7d1c/e2e/78b55644-0d76-42fb-a98a-4b0706e7c488.md", "", "", '230d5acb-504b-4637-9f8e-a5c285449e76.md')
$ws2.Hyperlinks.Add($ws2.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/4f3555da4b77e2e3cc6db45f26b53df474aa2aeb/e2e/78b55644-0d76-42fb-a98a-4b0706e7c488.md", "", "", '230d5acb-504b-4637-9f8e-a5c285449e76.md')
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d15cafe1ec021c40fea94da5b6f49b1a36cc7d1c/e2e/7b2c538e-3fdb-45ed-be82-d28ae0801cc8.md", "", "", 'ffffd46d8d01-c2cf-4e44-ab34-61cf8c3e76e4.md')
$ws2.Hyperlinks.Add($ws2.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/4f3555da4b77e2e3cc6db45f26b53df474aa2aeb/e2e/7b2c538e-3fdb-45ed-be82-d28ae0801cc8.md", "", "", 'ffffd46d8d01-c2cf-4e44-ab34-61cf8c3e76e4.md')

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Range("A2").Value = '230d5acb-504b-4637-9f8e-a5c285449e76.md'
$ws3.Range("G2").Value = '230d5acb-504b-4637-9f8e-a5c285449e76.f54858a962c9264f22fab170293849181df9a130.de-de.xlf'
$ws3.Range("H2").Value = '2016-09-04 07:09:09'
$ws3.Range("I2").Value = '230d5acb-504b-4637-9f8e-a5c285449e76.md'
$ws3.Range("J2").Value = '230d5acb-504b-4637-9f8e-a5c285449e76.f54858a962c9264f22fab170293849181df9a130.de-de.xlf'
$ws3.Range("K2").Value = '2016-09-04 07:09:36'

$ws3.Range("A3").Value = 'ffffd46d8d01-c2cf-4e44-ab34-61cf8c3e76e4.md'
$ws3.Range("G3").Value = '230d5acb-504b-4637-9f8e-a5c285449e76.f54858a962c9264f22fab170293849181df9a130.de-de.xlf'
$ws3.Range("H3").Value = '2016-09-04 07:09:09'
$ws3.Range("I3").Value = 'ffffd46d8d01-c2cf-4e44-ab34-61cf8c3e76e4.md'
$ws3.Range("J3").Value = '230d5acb-504b-4637-9f8e-a5c285449e76.f54858a962c9264f22fab170293849181df9a130.de-de.xlf'
$ws3.Range("K3").Value = '2016-09-04 07:09:36'

$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d15cafe1ec021c40fea94da5b6f49b1a36cc7d1c/e2e/78b55644-0d76-42fb-a98a-4b0706e7c488.md", "", "", '230d5acb-504b-4637-9f8e-a5c285449e76.md')
$ws3.Hyperlinks.Add($ws3.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/013323cbe01b0f576ad59859555a3f76c9593842/e2e/78b55644-0d76-42fb-a98a-4b0706e7c488.md", "", "", '230d5acb-504b-4637-9f8e-a5c285449e76.md')
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d15cafe1ec021c40fea94da5b6f49b1a36cc7d1c/e2e/7b2c538e-3fdb-45ed-be82-d28ae0801cc8.md", "", "", 'ffffd46d8d01-c2cf-4e44-ab34-61cf8c3e76e4.md')
$ws3.Hyperlinks.Add($ws3.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/013323cbe01b0f576ad59859555a3f76c9593842/e2e/7b2c538e-3fdb-45ed-be82-d28ae0801cc8.md", "", "", 'ffffd46d8d01-c2cf-4e44-ab34-61cf8c3e76e4.md')
